# Applies the "Proyecto 330" workbook update described by the commit diff:
#  - Sets a number of previously-blank payment cells (columns D/F/H,
#    "Cuota Enero/Febrero/Marzo") to 10000, reusing the existing currency
#    cell style (same style index as other populated cells in the sheet).
#  - Normalizes a few blank cells' style (D/F column, still blank) the same
#    way the source workbook does it (style index 5 -> 6).
#  - Appends two new people rows (357 "ZILLI / GRISELDA" and 358
#    "RODRIGUEZ / MABEL") at the bottom of the table, each with a 10000
#    payment in "Cuota Enero" and "Cuota Febrero".
#  - The totals (J column SUM formulas and the M2/O2 summary formulas) are
#    recalculated automatically by the engine once the source values change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell that already carries the "populated currency" style used
# throughout column D/F/H (style index 3 in the original workbook).
$populatedSrc = $ws.Range("D3")

# Reference cell that already carries the "blank, no currency format"
# style used for a few blank D/F/H cells (style index 6).
$blankSrc = $ws.Range("G2")

# Cells that go from blank to a 10000 payment (value + style copied from
# $populatedSrc, matching the existing look of every other populated cell).
$populateCells = @(
    "D10",
    "F18",
    "D36", "F36", "H36",
    "F56",
    "D70",
    "D71",
    "F76",
    "F78",
    "D93", "F93",
    "F100",
    "D116",
    "D172",
    "D173",
    "F174",
    "D177",
    "D247",
    "H279",
    "D318",
    "D339",
    "F343",
    "D352"
)

foreach ($addr in $populateCells) {
    $populatedSrc.Copy($ws.Range($addr))
}

# Cells that stay blank but change style (index 5 -> 6), same as the
# source diff.
$styleOnlyCells = @("F10", "F119", "F318")

foreach ($addr in $styleOnlyCells) {
    $blankSrc.Copy($ws.Range($addr))
}

# --- New row 357: ZILLI / GRISELDA ---------------------------------------
# Use row 341 (A/B/C style 2, D/F style 3) as a style template so the new
# rows pick up the exact same cell styles already used elsewhere.
$ws.Range("A341:C341").Copy($ws.Range("A357:C357"))
$ws.Range("D341").Copy($ws.Range("D357"))
$ws.Range("D341").Copy($ws.Range("F357"))

$ws.Range("A357").Value = 357
$ws.Range("B357").Value = "ZILLI "
$ws.Range("C357").Value = "GRISELDA"

# --- New row 358: RODRIGUEZ / MABEL --------------------------------------
$ws.Range("A341:C341").Copy($ws.Range("A358:C358"))
$ws.Range("D341").Copy($ws.Range("D358"))
$ws.Range("D341").Copy($ws.Range("F358"))
$ws.Range("B341").Copy($ws.Range("H358"))

$ws.Range("A358").Value = 357
$ws.Range("B358").Value = "RODRIGUEZ"
$ws.Range("C358").Value = "MABEL"
$ws.Range("H358").Value = " `$10.000,00"
